# SE-341 Validate the metabolomics data, not just the metadata.
#
# The "openbis-data" sheet previously held valid CHEBI compound identifiers
# in column A (rows 2-3). Replace them with the bogus placeholder "foo" so
# the metabolomics *data* rows exercise the same "bad data" validation path
# as the metadata sheet already does.

$wb = $excel.ActiveWorkbook

$metadataSheet = $wb.Worksheets.Item(1)   # "openbis-metadata"
$dataSheet     = $wb.Worksheets.Item(2)   # "openbis-data"

# Corrupt the CompoundID values used for the data rows.
$dataSheet.Range("A2").Value = "foo"
$dataSheet.Range("A3").Value = "foo"

# Move the selection/active view onto the data sheet, mirroring the
# reviewer now looking at (and validating) the data rows instead of the
# metadata sheet.
$metadataSheet.Range("D15").Select()

$dataSheet.Activate()
$dataSheet.Range("A4").Select()
